$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("n5", "n5_IMG_3177.jpeg", "'True", "no_meltpatch", "negative"),
    @("n6", "n6_IMG_3176.jpeg", "'True", "no_meltpatch", "negative"),
    @("n7", "n7_IMG_3179.jpeg", "'True", "no_meltpatch", "negative"),
    @("n8", "n8_IMG_3174.jpeg", "'True", "no_meltpatch", "negative")
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
